# Add a new bulleted item right after "...to add new user." that documents
# the new REST API (receives Userid/password, authenticates, responds to UI).
#
# We do this with Find/Replace rather than Range.InsertParagraphAfter()
# because this host's InsertParagraphAfter() mis-copies the bold run
# formatting that appears earlier in the source paragraph ("REST API").
# Using Find.Execute with a replacement string containing a paragraph mark
# (Chr(13)) splits the paragraph the same way a user pressing Enter then
# typing would, but keeps the newly typed run free of that stray formatting,
# and the new paragraph correctly inherits the surrounding ListParagraph /
# numPr (ilvl 0, numId 4) list formatting.

$d = $word.ActiveDocument

$wdReplaceOne = 1
$wdFindContinue = 1

$paragraphMark = [char]13
$newBullet = "Create REST API which will receive Userid and password from UI , authenticate and send back response to UI."
$replacement = " to add new user." + $paragraphMark + $newBullet

$d.Content.Find.Execute(
    " to add new user.",  # FindText
    $false,                # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    $wdFindContinue,       # Wrap
    $false,                # Format
    $replacement,          # ReplaceWith
    $wdReplaceOne          # Replace
) | Out-Null
